# lista_de_turnos.xlsx - "Actualización de la planilla de turnos"
#
# Several RUT (Chilean ID) values in column B were typed with a stray
# space before the final dash (e.g. "18.392.207 -6" instead of
# "18.392.207-6"). This corrects those entries so they match the
# already-correct RUT format used elsewhere in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$fixes = @{
    "B3"  = "18.392.207-6"
    "B20" = "16.759.697-5"
    "B23" = "16.751.516-9"
    "B37" = "16.759.697-5"
    "B43" = "16.751.516-9"
    "B49" = "18.392.207-6"
    "B54" = "16.759.697-5"
    "B58" = "16.751.516-9"
    "B64" = "18.392.207-6"
    "B68" = "16.759.697-5"
    "B74" = "16.751.516-9"
}

foreach ($addr in $fixes.Keys) {
    $ws.Range($addr).Value = $fixes[$addr]
}

# Reflect the author's on-screen state when the file was saved: zoomed in
# to 175% with cell F9 selected (instead of the previous K11 selection).
$null = $ws.Activate()
$excel.ActiveWindow.Zoom = 175
$null = $ws.Range("F9").Select()

$wb.Save()
